$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 8: update reporting period / emission dates ---
# B8: fecha de inicio del periodo (2022-01-01 -> 2022-04-01)
$ws.Range("B8").Value = 44652
# C8: fecha de termino del periodo (2022-03-31 -> 2022-06-30)
$ws.Range("C8").Value = 44742
# I8: fecha de validacion (2022-04-08 -> 2022-07-11)
$ws.Range("I8").Value = 44753
# J8: fecha de actualizacion (2022-04-08 -> 2022-07-11)
$ws.Range("J8").Value = 44753

# --- H8: area responsable text (drop leading space) ---
$ws.Range("H8").Value = "Secretaria Academica (UPP)"
# Style refresh Excel performs when the cell text is retyped/cleared of
# formatting: plain Calibri 11 black, default (general) horizontal align.
$ws.Range("H8").Font.Name = "Calibri"
$ws.Range("H8").Font.Size = 11
$ws.Range("H8").Font.Color = 0
$ws.Range("H8").HorizontalAlignment = 1

# --- K8: replace long note with the new shorter note ---
$ws.Range("K8").Value = "La Universidad Politécnica de Pachuca, no tiene opiniones y recomendaciones del Consejo Consultivo."

# --- Row heights ---
$ws.Rows.Item(3).RowHeight = 37.5
$ws.Rows.Item(8).RowHeight = 46.5

# --- Column K width ---
$ws.Columns.Item(11).ColumnWidth = 42.65

# --- Data validation: extend list validation down column D ---
$ws.Range("D8:D201").Validation.Delete()
$ws.Range("D8:D201").Validation.Add(3, 1, 1, "=Hidden_13")
$ws.Range("D8:D201").Validation.ShowInput = $false

# --- Selection / view state ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("K12").Select()
